$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Form the consolidated report: update the "Absent" column (H) values
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
